# Applies the cryptos-list price/volume refresh for this run.
# Column layout: A=index, B=Coin, C=Link, D=Price, E=Volume(1h)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price-column (D) values are digit-grouped strings (e.g. "29.387.77") that
# Excel would otherwise coerce to a number and mangle (dropping trailing
# zeroes, re-parsing the dots). Force literal text the same way a user would
# (leading apostrophe), then strip the resulting style back to Normal so the
# cell keeps the workbook default formatting.
function Set-TextValue($cell, $value) {
    $cell.Value = "'" + $value
    $cell.Style = "Normal"
}

# Row 2: D2, E2
Set-TextValue $ws.Cells.Item(2, 4) '29.387.77'
$ws.Cells.Item(2, 5).Value = '  +0.14%  '
# Row 3: D3, E3
Set-TextValue $ws.Cells.Item(3, 4) '1.841.75'
$ws.Cells.Item(3, 5).Value = '  -0.15%  '
# Row 4: E4
$ws.Cells.Item(4, 5).Value = '  +0.23%  '
# Row 5: D5, E5
Set-TextValue $ws.Cells.Item(5, 4) '239.17'
$ws.Cells.Item(5, 5).Value = '  -0.35%  '
# Row 6: D6, E6
Set-TextValue $ws.Cells.Item(6, 4) '0.6265'
$ws.Cells.Item(6, 5).Value = '  +0.03%  '
# Row 7: E7
$ws.Cells.Item(7, 5).Value = '  +0.19%  '
# Row 8: D8, E8
Set-TextValue $ws.Cells.Item(8, 4) '0.07423'
$ws.Cells.Item(8, 5).Value = '  -1.04%  '
# Row 9: D9, E9
Set-TextValue $ws.Cells.Item(9, 4) '0.2892'
$ws.Cells.Item(9, 5).Value = '  -0.23%  '
# Row 10: D10, E10
Set-TextValue $ws.Cells.Item(10, 4) '24.94'
$ws.Cells.Item(10, 5).Value = '  +2.02%  '
# Row 11: D11, E11
Set-TextValue $ws.Cells.Item(11, 4) '0.07719'
$ws.Cells.Item(11, 5).Value = '  -0.18%  '
# Row 12: D12, E12
Set-TextValue $ws.Cells.Item(12, 4) '1.831.20'
$ws.Cells.Item(12, 5).Value = '  -0.67%  '
# Row 13: D13, E13
Set-TextValue $ws.Cells.Item(13, 4) '4.976'
$ws.Cells.Item(13, 5).Value = '  -0.35%  '
# Row 14: D14, E14
Set-TextValue $ws.Cells.Item(14, 4) '0.6743'
$ws.Cells.Item(14, 5).Value = '  -1.00%  '
# Row 15: E15
$ws.Cells.Item(15, 5).Value = '  -2.39%  '
# Row 16: D16, E16
Set-TextValue $ws.Cells.Item(16, 4) '81.76'
$ws.Cells.Item(16, 5).Value = '  -0.36%  '
# Row 17: D17, E17
Set-TextValue $ws.Cells.Item(17, 4) '6.212'
$ws.Cells.Item(17, 5).Value = '  +0.36%  '
# Row 18: D18, E18
Set-TextValue $ws.Cells.Item(18, 4) '29.421.49'
$ws.Cells.Item(18, 5).Value = '  +0.12%  '
# Row 19: D19, E19
Set-TextValue $ws.Cells.Item(19, 4) '234.03'
$ws.Cells.Item(19, 5).Value = '  +2.19%  '
# Row 20: D20, E20
Set-TextValue $ws.Cells.Item(20, 4) '12.31'
$ws.Cells.Item(20, 5).Value = '  -0.16%  '
# Row 21: E21
$ws.Cells.Item(21, 5).Value = '  +0.20%  '
# Row 22: D22, E22
Set-TextValue $ws.Cells.Item(22, 4) '7.291'
$ws.Cells.Item(22, 5).Value = '  -2.73%  '
# Row 23: E23
$ws.Cells.Item(23, 5).Value = '  +0.17%  '
# Row 24: D24, E24
Set-TextValue $ws.Cells.Item(24, 4) '158.60'
$ws.Cells.Item(24, 5).Value = '  -0.07%  '
# Row 25: D25, E25
Set-TextValue $ws.Cells.Item(25, 4) '8.494'
$ws.Cells.Item(25, 5).Value = '  +0.79%  '
# Row 26: D26
Set-TextValue $ws.Cells.Item(26, 4) '0.1345'
# Row 27: E27
$ws.Cells.Item(27, 5).Value = '  -1.22%  '
# Row 28: D28, E28
Set-TextValue $ws.Cells.Item(28, 4) '0.07294'
$ws.Cells.Item(28, 5).Value = '  +12.86%  '
# Row 29: D29, E29
Set-TextValue $ws.Cells.Item(29, 4) '1.466'
$ws.Cells.Item(29, 5).Value = '  +4.39%  '
# Row 30: E30
$ws.Cells.Item(30, 5).Value = '  -0.05%  '
# Row 31: D31, E31
Set-TextValue $ws.Cells.Item(31, 4) '4.040'
$ws.Cells.Item(31, 5).Value = '  -1.33%  '
# Row 32: D32
Set-TextValue $ws.Cells.Item(32, 4) '4.029'
# Row 33: D33, E33
Set-TextValue $ws.Cells.Item(33, 4) '1.815'
$ws.Cells.Item(33, 5).Value = '  -0.94%  '
# Row 34: E34
$ws.Cells.Item(34, 5).Value = '  -0.11%  '
# Row 35: D35, E35
Set-TextValue $ws.Cells.Item(35, 4) '0.6972'
$ws.Cells.Item(35, 5).Value = '  -0.11%  '
# Row 36: D36
Set-TextValue $ws.Cells.Item(36, 4) '2.571'
# Row 37: E37
$ws.Cells.Item(37, 5).Value = '  +0.35%  '
# Row 38: D38, E38
Set-TextValue $ws.Cells.Item(38, 4) '6.912'
$ws.Cells.Item(38, 5).Value = '  +2.26%  '
# Row 39: E39
$ws.Cells.Item(39, 5).Value = '  -0.87%  '
# Row 40: D40, E40
Set-TextValue $ws.Cells.Item(40, 4) '1.232.21'
$ws.Cells.Item(40, 5).Value = '  -2.86%  '
# Row 41: D41, E41
Set-TextValue $ws.Cells.Item(41, 4) '0.9597'
$ws.Cells.Item(41, 5).Value = '  +4.84%  '
# Row 42: E42
$ws.Cells.Item(42, 5).Value = '  +0.22%  '
# Row 43: D43, E43
Set-TextValue $ws.Cells.Item(43, 4) '1.997.14'
$ws.Cells.Item(43, 5).Value = '  -0.43%  '
# Row 44: D44, E44
Set-TextValue $ws.Cells.Item(44, 4) '100.93'
$ws.Cells.Item(44, 5).Value = '  -0.36%  '
# Row 45: D45, E45
Set-TextValue $ws.Cells.Item(45, 4) '65.43'
$ws.Cells.Item(45, 5).Value = '  -1.25%  '
# Row 46: B46, C46, D46, E46
$ws.Cells.Item(46, 2).Value = 'BabyDogeCoin'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextValue $ws.Cells.Item(46, 4) '0.00000000120'
$ws.Cells.Item(46, 5).Value = '  +7.95%  '
# Row 47: B47, C47, D47, E47
$ws.Cells.Item(47, 2).Value = 'RenderToken'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws.Cells.Item(47, 4) '1.715'
$ws.Cells.Item(47, 5).Value = '  -0.76%  '
# Row 48: B48, C48, D48, E48
$ws.Cells.Item(48, 2).Value = 'Aptos'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue $ws.Cells.Item(48, 4) '6.950'
$ws.Cells.Item(48, 5).Value = '  -1.87%  '
# Row 49: E49
$ws.Cells.Item(49, 5).Value = '  -2.70%  '
# Row 50: B50, C50, D50, E50
$ws.Cells.Item(50, 2).Value = 'EnergySwap'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Cells.Item(50, 4) '8.851'
$ws.Cells.Item(50, 5).Value = '  -1.38%  '
# Row 51: B51, C51, D51, E51
$ws.Cells.Item(51, 2).Value = 'TheSandbox'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue $ws.Cells.Item(51, 4) '0.3897'
$ws.Cells.Item(51, 5).Value = '  -1.73%  '
